$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 1868
$ws.Range("J112").Value = 2085
$ws.Range("L112").Value = 6255
$ws.Range("N112").Value = -8471

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 16390
$ws.Range("J24").Value = 16390
$ws.Range("L24").Value = 16390
$ws.Range("N24").Value = -17138
# Row 32
$ws.Range("H32").Value = 1395246.4
$ws.Range("I32").Value = 1547873.5
$ws.Range("J32").Value = 21602.166
$ws.Range("K32").Value = 1547873.5
$ws.Range("L32").Value = 21602.166
$ws.Range("M32").Value = -1547586.5
$ws.Range("N32").Value = -22176.166
# Row 88
$ws.Range("H88").Value = 16901.568
$ws.Range("I88").Value = 20710.629
$ws.Range("J88").Value = 2088.5557
$ws.Range("K88").Value = 20710.629
$ws.Range("L88").Value = 2088.5557
$ws.Range("M88").Value = -20304.629
$ws.Range("N88").Value = -2900.5557
# Row 91
$ws.Range("H91").Value = 16901.568
$ws.Range("I91").Value = 20710.629
$ws.Range("J91").Value = 2088.5557
$ws.Range("K91").Value = 20710.629
$ws.Range("L91").Value = 2088.5557
$ws.Range("M91").Value = -19306.629
$ws.Range("N91").Value = -4896.5557
# Row 100
$ws.Range("H100").Value = 16390
$ws.Range("J100").Value = 16390
$ws.Range("L100").Value = 16390
$ws.Range("N100").Value = -18554
# Row 113
$ws.Range("H113").Value = 39330
$ws.Range("J113").Value = 39330
$ws.Range("L113").Value = 39330
$ws.Range("N113").Value = -48008
# Row 132
$ws.Range("H132").Value = 25717.166
$ws.Range("I132").Value = 35708.484
$ws.Range("K132").Value = 107125.452
$ws.Range("M132").Value = -104595.452

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1024.8846
$ws.Range("I20").Value = 802.93335
$ws.Range("J20").Value = 1327.5454
$ws.Range("K20").Value = 802.93335
$ws.Range("L20").Value = 1327.5454
$ws.Range("M20").Value = -555.93335
$ws.Range("N20").Value = -1821.5454

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1035576.1
$ws.Range("I31").Value = 1043.5
$ws.Range("J31").Value = 1570679.2
$ws.Range("K31").Value = 1043.5
$ws.Range("L31").Value = 1570679.2
$ws.Range("M31").Value = -748.5
$ws.Range("N31").Value = -1571269.2
# Row 34
$ws.Range("H34").Value = 1035576.1
$ws.Range("I34").Value = 1043.5
$ws.Range("J34").Value = 1570679.2
$ws.Range("K34").Value = 1043.5
$ws.Range("L34").Value = 1570679.2
$ws.Range("M34").Value = -841.5
$ws.Range("N34").Value = -1571083.2
# Row 94
$ws.Range("H94").Value = 4338
$ws.Range("I94").Value = 977.75
$ws.Range("J94").Value = 6018.125
$ws.Range("K94").Value = 977.75
$ws.Range("L94").Value = 6018.125
$ws.Range("M94").Value = -526.75
$ws.Range("N94").Value = -6920.125

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 31
$ws.Range("H31").Value = 707.2727
$ws.Range("I31").Value = 707.2727
$ws.Range("K31").Value = 2121.8181
$ws.Range("M31").Value = -1833.8181
# Row 49
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 2000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 6000
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -6312
# Row 57
$ws.Range("H57").Value = 2928.5715
$ws.Range("I57").Value = 500
$ws.Range("J57").Value = 4750
$ws.Range("K57").Value = 1500
$ws.Range("L57").Value = 14250
$ws.Range("M57").Value = -941
$ws.Range("N57").Value = -15368
# Row 74
$ws.Range("H74").Value = 5999.6665
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 5999.6665
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 17998.9995
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -20120.9995
# Row 77
$ws.Range("H77").Value = 5999.6665
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 5999.6665
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 53996.9985
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -64604.9985
# Row 96
$ws.Range("H96").Value = 5000
$ws.Range("J96").Value = 5000
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -19118
# Row 106
$ws.Range("H106").Value = 4077.5715
$ws.Range("I106").Value = 2633.3333
$ws.Range("J106").Value = 4318.278
$ws.Range("K106").Value = 7899.999899999999
$ws.Range("L106").Value = 12954.834
$ws.Range("M106").Value = -6953.999899999999
$ws.Range("N106").Value = -14846.834
# Row 107
$ws.Range("H107").Value = 1033.6875
$ws.Range("J107").Value = 2214.3635
$ws.Range("L107").Value = 6643.0905
$ws.Range("N107").Value = -10483.0905
# Row 130
$ws.Range("H130").Value = 126039.125
$ws.Range("I130").Value = 606
$ws.Range("J130").Value = 335094.34
$ws.Range("K130").Value = 1818
$ws.Range("L130").Value = 1005283.02
$ws.Range("M130").Value = 3202
$ws.Range("N130").Value = -1015323.02
# Row 131
$ws.Range("H131").Value = 1231.5
$ws.Range("I131").Value = 1189.3636
$ws.Range("J131").Value = 1241.8
$ws.Range("K131").Value = 3568.0908
$ws.Range("L131").Value = 3725.4
$ws.Range("M131").Value = 1471.9092
$ws.Range("N131").Value = -13805.4

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 18
$ws.Range("H18").Value = 5200.8335
$ws.Range("I18").Value = 5200.8335
$ws.Range("K18").Value = 5200.8335
$ws.Range("M18").Value = -5028.8335
# Row 20
$ws.Range("H20").Value = 3666.6667
$ws.Range("I20").Value = 3000
$ws.Range("J20").Value = 4000
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 4000
$ws.Range("M20").Value = -2774
$ws.Range("N20").Value = -4452
# Row 22
$ws.Range("H22").Value = 427
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 357.8
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 357.8
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -947.8
# Row 27
$ws.Range("H27").Value = 427
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 357.8
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 357.8
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -571.8
# Row 36
$ws.Range("H36").Value = 36991.5
$ws.Range("J36").Value = 36991.5
$ws.Range("L36").Value = 36991.5
$ws.Range("N36").Value = -38115.5
# Row 132
$ws.Range("H132").Value = 9016922
$ws.Range("I132").Value = 2599.3635
$ws.Range("J132").Value = 22237930
$ws.Range("K132").Value = 7798.0905
$ws.Range("L132").Value = 66713790
$ws.Range("M132").Value = -5268.0905
$ws.Range("N132").Value = -66718850
